# Update "handback status" timestamps to reflect a freshly generated report.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date"
$overview.Range("G2").Value = "2016-08-22 21:05:25"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$zhcn.Range("H2").Value = "2016-08-22 21:05:18"
$zhcn.Range("K2").Value = "2016-08-22 21:05:46"

# de-de sheet: "Correspond Handoff Datetime" (shares the same value as Overview!G2)
# and "Correspond Handback DateTime"
$dede.Range("H2").Value = "2016-08-22 21:05:25"
$dede.Range("K2").Value = "2016-08-22 21:05:54"
